$wb = $excel.ActiveWorkbook

# --- test_suite sheet (sheet1) ---
$ws1 = $wb.Worksheets.Item("test_suite")
$ws1.Range("A5").Value = "Export_Daily_Reports"
$ws1.Range("B7").Value = "N"
$ws1.Range("B8").Value = "N"
$ws1.Range("B14").Value = "N"

# --- Reports sheet (sheet2) ---
$ws2 = $wb.Worksheets.Item("Reports")
$ws2.Range("B2").Value = "JANUARY"
$ws2.Range("B2").Select()

# test_suite becomes the active sheet/tab, with B6 selected
$ws1.Activate()
$ws1.Range("B6").Select()
